$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J5").Value = "."
$ws.Range("J5").Select() | Out-Null
